# Update "想去人数" (want-to-go count) figures in column F across the
# "展览", "演出" and "全部类型" sheets, matching the refreshed data pull.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 555
$ws1.Range("F9").Value  = 1438
$ws1.Range("F12").Value = 2993
$ws1.Range("F13").Value = 391
$ws1.Range("F14").Value = 1605
$ws1.Range("F15").Value = 1358
$ws1.Range("F16").Value = 786
$ws1.Range("F17").Value = 235
$ws1.Range("F18").Value = 1368
$ws1.Range("F22").Value = 396
$ws1.Range("F23").Value = 3459
$ws1.Range("F24").Value = 676

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 78

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F18").Value = 555
$ws4.Range("F19").Value = 1438
$ws4.Range("F22").Value = 2993
$ws4.Range("F23").Value = 391
$ws4.Range("F24").Value = 1605
$ws4.Range("F25").Value = 1358
$ws4.Range("F26").Value = 786
$ws4.Range("F27").Value = 235
$ws4.Range("F28").Value = 1368
$ws4.Range("F34").Value = 396
$ws4.Range("F35").Value = 3459
$ws4.Range("F36").Value = 676
$ws4.Range("F39").Value = 78
